$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, matching the formatting of the existing
# header cells (bold, centered, bordered) by copying G1's format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill H2:H12 with 0 (the new Save values, unstyled like the other
# numeric columns).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
